$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("M12").ClearContents()
$ws.Range("N12").ClearContents()
$ws.Range("H38").Value = 2365
$ws.Range("I38").Value = 59.375
$ws.Range("J38").Value = 5000
$ws.Range("K38").Value = 178.125
$ws.Range("L38").Value = 15000
$ws.Range("M38").Value = 193.875
$ws.Range("N38").Value = -15744
$ws.Range("H43").Value = 46314.824
$ws.Range("J43").Value = 127558.664
$ws.Range("L43").Value = 127558.664
$ws.Range("N43").Value = -127696.664
$ws.Range("H51").Value = 0
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("M51").ClearContents()
$ws.Range("N51").ClearContents()
$ws.Range("H98").Value = 1817.6666
$ws.Range("I98").Value = 1525.2858
$ws.Range("K98").Value = 1525.2858
$ws.Range("M98").Value = -27.28580000000011
$ws.Range("H106").Value = 4987.25
$ws.Range("I106").Value = 4987.25
$ws.Range("K106").Value = 4987.25
$ws.Range("M106").Value = -4356.25
$ws.Range("H113").Value = 4442.7
$ws.Range("I113").Value = 4179.6
$ws.Range("J113").Value = 4705.8
$ws.Range("K113").Value = 4179.6
$ws.Range("L113").Value = 4705.8
$ws.Range("M113").Value = -925.6000000000004
$ws.Range("N113").Value = -11213.8
$ws.Range("H122").Value = 1817.6666
$ws.Range("I122").Value = 1525.2858
$ws.Range("K122").Value = 4575.857400000001
$ws.Range("M122").Value = -2125.857400000001
$ws.Range("H132").Value = 5296
$ws.Range("I132").Value = 5196.857
$ws.Range("K132").Value = 15590.571
$ws.Range("M132").Value = -13060.571
$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("M136").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 31
$ws.Range("J5").Value = 32.333332
$ws.Range("L5").Value = 32.333332
$ws.Range("N5").Value = -256.333332
$ws.Range("H32").Value = 9679.467000000001
$ws.Range("J32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("N32").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 31
$ws.Range("J4").Value = 32.333332
$ws.Range("L4").Value = 32.333332
$ws.Range("N4").Value = -262.333332
$ws.Range("H86").Value = 1928.5652
$ws.Range("J86").Value = 3910.3333
$ws.Range("L86").Value = 3910.3333
$ws.Range("N86").Value = -6156.3333
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()
$ws.Range("H89").Value = 1928.5652
$ws.Range("J89").Value = 3910.3333
$ws.Range("L89").Value = 19551.6665
$ws.Range("N89").Value = -30783.6665
$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()
$ws.Range("H99").Value = 2697.2727
$ws.Range("I99").Value = 2630
$ws.Range("K99").Value = 2630
$ws.Range("M99").Value = -1132
$ws.Range("H105").Value = 1839.4
$ws.Range("I105").Value = 1774.5
$ws.Range("J105").Value = 1882.6666
$ws.Range("K105").Value = 1774.5
$ws.Range("L105").Value = 1882.6666
$ws.Range("M105").Value = -27.5
$ws.Range("N105").Value = -5376.6666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H74").Value = 34999.9
$ws.Range("J74").Value = 34999.9
$ws.Range("L74").Value = 34999.9
$ws.Range("N74").Value = -36747.9
$ws.Range("H77").Value = 34999.9
$ws.Range("J77").Value = 34999.9
$ws.Range("L77").Value = 104999.7
$ws.Range("N77").Value = -113735.7
$ws.Range("H103").Value = 4499.5
$ws.Range("I103").Value = 4499.5
$ws.Range("K103").Value = 4499.5
$ws.Range("M103").Value = -3327.5
$ws.Range("H108").Value = 49995
$ws.Range("I108").Value = 20000
$ws.Range("J108").Value = 79990
$ws.Range("K108").Value = 20000
$ws.Range("L108").Value = 79990
$ws.Range("M108").Value = -16160
$ws.Range("N108").Value = -87670
$ws.Range("H121").Value = 47184.332
$ws.Range("J121").Value = 50777
$ws.Range("L121").Value = 50777
$ws.Range("N121").Value = -53397
$ws.Range("H122").Value = 1043.6666
$ws.Range("I122").Value = 727.5714
$ws.Range("K122").Value = 2182.7142
$ws.Range("M122").Value = 267.2857999999997
$ws.Range("H132").Value = 1450
$ws.Range("J132").Value = 1400
$ws.Range("L132").Value = 4200
$ws.Range("N132").Value = -9260
$ws.Range("H134").Value = 11376.5
$ws.Range("I134").Value = 10169.5
$ws.Range("J134").Value = 14997.5
$ws.Range("K134").Value = 30508.5
$ws.Range("L134").Value = 44992.5
$ws.Range("M134").Value = -27973.5
$ws.Range("N134").Value = -50062.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H33").Value = 10999
$ws.Range("I33").Value = 0
$ws.Range("J33").Value = 10999
$ws.Range("K33").Value = 0
$ws.Range("L33").Value = 10999
$ws.Range("M33").ClearContents()
$ws.Range("N33").Value = -11503
$ws.Range("H113").Value = 1933.3334
$ws.Range("I113").Value = 1800
$ws.Range("K113").Value = 1800
$ws.Range("M113").Value = 370
$ws.Range("H132").Value = 1960
$ws.Range("I132").Value = 1960
$ws.Range("K132").Value = 5880
$ws.Range("M132").Value = -3350

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2333.1667
$ws.Range("J22").Value = 3999
$ws.Range("L22").Value = 3999
$ws.Range("N22").Value = -4589
$ws.Range("H27").Value = 2333.1667
$ws.Range("J27").Value = 3999
$ws.Range("L27").Value = 3999
$ws.Range("N27").Value = -4213
$ws.Range("H46").Value = 3230.7693
$ws.Range("I46").Value = 3000
$ws.Range("J46").Value = 4500
$ws.Range("K46").Value = 3000
$ws.Range("L46").Value = 4500
$ws.Range("M46").Value = -2812
$ws.Range("N46").Value = -4876
$ws.Range("H61").Value = 4152.75
$ws.Range("I61").Value = 4288.7856
$ws.Range("K61").Value = 4288.7856
$ws.Range("M61").Value = -4086.7856
$ws.Range("H87").Value = 25000
$ws.Range("J87").Value = 25000
$ws.Range("L87").Value = 25000
$ws.Range("N87").Value = -27246
$ws.Range("H90").Value = 25000
$ws.Range("J90").Value = 25000
$ws.Range("L90").Value = 75000
$ws.Range("N90").Value = -86232
$ws.Range("H93").Value = 417.6
$ws.Range("I93").Value = 417.6
$ws.Range("K93").Value = 417.6
$ws.Range("M93").Value = 830.4
$ws.Range("H113").Value = 4152.75
$ws.Range("I113").Value = 4288.7856
$ws.Range("K113").Value = 4288.7856
$ws.Range("M113").Value = -2118.7856

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 691.5
$ws.Range("I81").Value = 549.6667
$ws.Range("K81").Value = 1099.3334
$ws.Range("M81").Value = -38.33339999999998
$ws.Range("H84").Value = 691.5
$ws.Range("I84").Value = 549.6667
$ws.Range("K84").Value = 5496.666999999999
$ws.Range("M84").Value = -192.6669999999995
$ws.Range("H122").Value = 2035.5454
$ws.Range("J122").Value = 1200
$ws.Range("L122").Value = 3600
$ws.Range("N122").Value = -8500

